$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.519.98"
$ws.Range("E2").Value = "  -5.02%  "

$ws.Range("D3").Value = "2.650.48"
$ws.Range("E3").Value = "  +1.53%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").Value = "'305.92"
$ws.Range("E5").Value = "  -0.41%  "

$ws.Range("D6").Value = "'96.75"
$ws.Range("E6").Value = "  -3.91%  "

$ws.Range("E7").Value = "  -2.33%  "

$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").Value = "'0.567"
$ws.Range("E9").Value = "  -2.65%  "

$ws.Range("D10").Value = "'37.58"
$ws.Range("E10").Value = "  -4.62%  "

$ws.Range("E11").Value = "  -2.16%  "

$ws.Range("D12").Value = "'7.93"
$ws.Range("E12").Value = "  -3.17%  "

$ws.Range("D13").Value = "3.058.94"
$ws.Range("E13").Value = "  +1.58%  "

$ws.Range("E14").Value = "  +0.75%  "

$ws.Range("D15").Value = "2.658.60"
$ws.Range("E15").Value = "  +1.70%  "

$ws.Range("D16").Value = "'0.909"
$ws.Range("E16").Value = "  -1.16%  "

$ws.Range("D17").Value = "'14.82"
$ws.Range("E17").Value = "  -1.20%  "

$ws.Range("D18").Value = "44.532.98"
$ws.Range("E18").Value = "  -5.35%  "

$ws.Range("D19").Value = "'6.78"
$ws.Range("E19").Value = "  +1.66%  "

$ws.Range("D20").Value = "0.0₃0994"
$ws.Range("E20").Value = "  -1.75%  "

$ws.Range("D21").Value = "'12.46"
$ws.Range("E21").Value = "  -3.55%  "

$ws.Range("D22").Value = "'74.19"
$ws.Range("E22").Value = "  +2.88%  "

$ws.Range("D23").Value = "'274.56"
$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("E24").Value = "  +5.03%  "

$ws.Range("D25").Value = "'2.99"
$ws.Range("E25").Value = "  -1.00%  "

$ws.Range("D26").Value = "'30.46"
$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").Value = "'10.37"
$ws.Range("E28").Value = "  -1.45%  "

$ws.Range("E29").Value = "  -2.89%  "

$ws.Range("D30").Value = "'37.40"
$ws.Range("E30").Value = "  -3.79%  "

$ws.Range("E31").Value = "  -0.36%  "

$ws.Range("D32").Value = "'3.71"
$ws.Range("E32").Value = "  +2.52%  "

$ws.Range("D33").Value = "'2.30"
$ws.Range("E33").Value = "  +5.37%  "

$ws.Range("D34").Value = "'153.26"
$ws.Range("E34").Value = "  +1.96%  "

$ws.Range("D35").Value = "'2.81"
$ws.Range("E35").Value = "  -2.15%  "

$ws.Range("D36").Value = "'0.0825"
$ws.Range("E36").Value = "  -2.02%  "

$ws.Range("E37").Value = "  -5.56%  "

$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").Value = "'24.99"
$ws.Range("E38").Value = "  +13.87%  "

$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.122"
$ws.Range("E39").Value = "  -0.32%  "

$ws.Range("D40").Value = "'15.76"
$ws.Range("E40").Value = "  +0.53%  "

$ws.Range("D41").Value = "'3.54"
$ws.Range("E41").Value = "  -2.13%  "

$ws.Range("E42").Value = "  -3.34%  "

$ws.Range("D43").Value = "2.125.57"
$ws.Range("E43").Value = "  -1.80%  "

$ws.Range("E44").Value = "  -6.07%  "

$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  -0.08%  "

$ws.Range("D46").Value = "'91.39"
$ws.Range("E46").Value = "  -4.62%  "

$ws.Range("D47").Value = "'9.31"
$ws.Range("E47").Value = "  -4.28%  "

$ws.Range("D48").Value = "2.910.34"
$ws.Range("E48").Value = "  +1.74%  "

$ws.Range("D49").Value = "'109.30"
$ws.Range("E49").Value = "  +0.50%  "

$ws.Range("D50").Value = "'1.59"
$ws.Range("E50").Value = "  -0.65%  "

$ws.Range("D51").Value = "'0.195"
$ws.Range("E51").Value = "  -2.02%  "

